$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "UC09 - Revender Produtos Recondicionados" -> split the trailing run so
#    the word "Produtos" becomes three runs: "Produto" / "s" / " Recondicionados"
#    (same text overall, just broken into separate runs - mirrors the diff).
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Produtos Recondicionados", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $rng.Start
    # Force a run split around the "s" in "Produtos" by toggling a character
    # property off then back on for just that letter - the host splits the
    # surrounding run into independent runs with identical formatting.
    $sChar = $d.Range($start + 7, $start + 8)
    $sChar.Font.Bold = 0
    $sChar.Font.Bold = 1
}

# ---------------------------------------------------------------------------
# 2) Resumo row: justify both cells (label + value) - add <w:jc w:val="both"/>
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
$resumoRow = $null
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $labelText = $t.Rows.Item($i).Cells.Item(1).Range.Text
    if ($labelText -like "Resumo*") {
        $resumoRow = $t.Rows.Item($i)
        break
    }
}
if ($resumoRow -ne $null) {
    $resumoRow.Cells.Item(1).Range.Paragraphs.Item(1).Alignment = 3
    $resumoRow.Cells.Item(2).Range.Paragraphs.Item(1).Alignment = 3
}

# ---------------------------------------------------------------------------
# 3) Merge "Sistema retorna " + bookmark(_GoBack) + "mensagem dizendo que o
#    produto sera analisado." into a single run / remove the _GoBack bookmark.
# ---------------------------------------------------------------------------
$target = "Sistema retorna mensagem dizendo que o produto será analisado."
$rng2 = $d.Content
$found2 = $rng2.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, $target, 2)
